# Fruta / hortaliza, semanal
# Re-order the daily price rows (2-26) according to the new weekly source
# extract. Only columns D, I, J, K, L, M, O, P change per row; every other
# column (A, B, C, E, F, G, H, N, Q, R) is constant across all rows in this
# sheet, so this is effectively a row permutation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: new row number -> array of values for columns D, I, J, K, L, M, O, P
$rowData = @{
    2  = @(44308, "Primera", 75,  5000, 5000, 5000, "Región del Maule",     312)
    3  = @(44398, "Primera", 80,  7000, 7000, 7000, "Región Metropolitana", 438)
    4  = @(44348, "Primera", 35,  7000, 7000, 7000, "Región del Maule",     438)
    5  = @(44315, "Primera", 40,  7000, 7000, 7000, "Región del Maule",     438)
    6  = @(44313, "Primera", 20,  7000, 7000, 7000, "Región del Maule",     438)
    7  = @(44354, "Primera", 100, 8000, 9000, 8500, "Región Metropolitana", 531)
    8  = @(44354, "Primera", 80,  9000, 9000, 9000, "Región del Maule",     562)
    9  = @(44314, "Segunda", 20,  5000, 5000, 5000, "Región del Maule",     312)
    10 = @(44392, "Primera", 95,  7000, 7000, 7000, "Región del Maule",     438)
    11 = @(44369, "Primera", 60,  7000, 7000, 7000, "Región Metropolitana", 438)
    12 = @(44399, "Primera", 80,  7000, 7000, 7000, "Región Metropolitana", 438)
    13 = @(44397, "Primera", 40,  8000, 8000, 8000, "Región Metropolitana", 500)
    14 = @(44420, "Primera", 45,  8000, 8000, 8000, "Región Metropolitana", 500)
    15 = @(44403, "Primera", 35,  5000, 5000, 5000, "Región Metropolitana", 312)
    16 = @(44312, "Primera", 40,  7000, 7000, 7000, "Región del Maule",     438)
    17 = @(44362, "Primera", 25,  8000, 8000, 8000, "Región Metropolitana", 500)
    18 = @(44305, "Primera", 35,  7000, 7000, 7000, "Región del Maule",     438)
    19 = @(44355, "Primera", 30,  8000, 8000, 8000, "Región Metropolitana", 500)
    20 = @(44372, "Primera", 50,  6000, 7000, 6400, "Región del Maule",     400)
    21 = @(44385, "Primera", 100, 7000, 7000, 7000, "Región del Maule",     438)
    22 = @(44371, "Primera", 200, 7000, 7000, 7000, "Región Metropolitana", 438)
    23 = @(44386, "Primera", 40,  7000, 7000, 7000, "Región del Maule",     438)
    24 = @(44389, "Primera", 55,  7000, 7000, 7000, "Región del Maule",     438)
    25 = @(44396, "Primera", 80,  7000, 7000, 7000, "Región Metropolitana", 438)
    26 = @(44467, "Primera", 40,  7000, 7000, 7000, "Región del Maule",     438)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 9).Value  = $vals[1]   # I - Calidad
    $ws.Cells.Item($row, 10).Value = $vals[2]   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[3]   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[4]   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $vals[6]   # O - Origen
    $ws.Cells.Item($row, 16).Value = $vals[7]   # P - Precio $/Kg
}
